# Update "analyse des variables (conservation/exclusion)" dashboard rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Outils column (G) gets a new tool entry "Overleaf" for the data-description rows
$ws.Range("G4").Value = "Overleaf"
$ws.Range("G5").Value = "Overleaf"
$ws.Range("G6").Value = "Overleaf"

# Row 7 (Préparation des données): status moves from "X" to "EC" and a start date is set
$ws.Range("D7").Value = "EC"
$ws.Range("E7").Value = (Get-Date -Year 2024 -Month 9 -Day 30 -Hour 0 -Minute 0 -Second 0)

# Row 8 (Exploration): tool name updated to include Word
$ws.Range("G8").Value = "Python / Anaconda / VSC / Word"

# Row 9 (Nettoyage): responsible changed and tool name updated to include Word
$ws.Range("C9").Value = "Alexis / Walid"
$ws.Range("G9").Value = "Python / Anaconda / VSC / Word"

# Row 10 (Fusion): responsible assigned, status moves from "X" to "EC", start date set
$ws.Range("C10").Value = "Abdelghani"
$ws.Range("D10").Value = "EC"
$ws.Range("E10").Value = (Get-Date -Year 2024 -Month 10 -Day 6 -Hour 0 -Minute 0 -Second 0)

# Leave the active selection on E18 to match the saved view state
$ws.Range("E18").Select()
